$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2106.7778
$ws.Range("I40").Value = 2196.25
$ws.Range("J40").Value = 2035.2
$ws.Range("K40").Value = 2196.25
$ws.Range("L40").Value = 2035.2
$ws.Range("M40").Value = -2021.25
$ws.Range("N40").Value = -2385.2
$ws.Range("H130").Value = 65279.168
$ws.Range("J130").Value = 65279.168
$ws.Range("L130").Value = 65279.168
$ws.Range("N130").Value = -75319.16800000001
$ws.Range("H131").Value = 3051.182
$ws.Range("I131").Value = 2179.7778
$ws.Range("J131").Value = 6972.5
$ws.Range("K131").Value = 6539.3334
$ws.Range("L131").Value = 20917.5
$ws.Range("M131").Value = -1499.3334
$ws.Range("N131").Value = -30997.5
$ws.Range("H138").Value = 9651868
$ws.Range("J138").Value = 13891595
$ws.Range("L138").Value = 41674785
$ws.Range("N138").Value = -41685065
$ws.Range("H141").Value = 4391.5835
$ws.Range("I141").Value = 2957.5293
$ws.Range("J141").Value = 7874.2856
$ws.Range("K141").Value = 8872.5879
$ws.Range("L141").Value = 23622.8568
$ws.Range("M141").Value = -3692.5879
$ws.Range("N141").Value = -33982.8568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490
$ws.Range("H132").Value = 3391.4583
$ws.Range("I132").Value = 2969.25
$ws.Range("K132").Value = 8907.75
$ws.Range("M132").Value = -6377.75
$ws.Range("H139").Value = 49267.4
$ws.Range("J139").Value = 49267.4
$ws.Range("L139").Value = 49267.4
$ws.Range("N139").Value = -59547.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1081.76
$ws.Range("I20").Value = 883.5
$ws.Range("J20").Value = 1434.2222
$ws.Range("K20").Value = 883.5
$ws.Range("L20").Value = 1434.2222
$ws.Range("M20").Value = -636.5
$ws.Range("N20").Value = -1928.2222
$ws.Range("H98").Value = 2492
$ws.Range("I98").Value = 2492
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2492
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("M98").Value = 503

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5636.8286
$ws.Range("I31").Value = 2064.55
$ws.Range("J31").Value = 10399.866
$ws.Range("K31").Value = 2064.55
$ws.Range("L31").Value = 10399.866
$ws.Range("M31").Value = -1769.55
$ws.Range("N31").Value = -10989.866
$ws.Range("H34").Value = 5636.8286
$ws.Range("I34").Value = 2064.55
$ws.Range("J34").Value = 10399.866
$ws.Range("K34").Value = 2064.55
$ws.Range("L34").Value = 10399.866
$ws.Range("M34").Value = -1862.55
$ws.Range("N34").Value = -10803.866
$ws.Range("H58").Value = 2715.577
$ws.Range("I58").Value = 1738.3846
$ws.Range("J58").Value = 3692.7693
$ws.Range("K58").Value = 1738.3846
$ws.Range("L58").Value = 3692.7693
$ws.Range("M58").Value = -1535.3846
$ws.Range("N58").Value = -4098.7693
$ws.Range("H132").Value = 3395.111
$ws.Range("I132").Value = 3171.4482
$ws.Range("K132").Value = 9514.3446
$ws.Range("M132").Value = -6984.3446
$ws.Range("H136").Value = 2715.577
$ws.Range("I136").Value = 1738.3846
$ws.Range("J136").Value = 3692.7693
$ws.Range("K136").Value = 5215.1538
$ws.Range("L136").Value = 11078.3079
$ws.Range("M136").Value = -2665.1538
$ws.Range("N136").Value = -16178.3079

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 17672.143
$ws.Range("J68").Value = 425.75
$ws.Range("L68").Value = 1277.25
$ws.Range("N68").Value = -2899.25
$ws.Range("H71").Value = 17672.143
$ws.Range("J71").Value = 425.75
$ws.Range("L71").Value = 3831.75
$ws.Range("N71").Value = -11943.75
$ws.Range("H131").Value = 8548524
$ws.Range("I131").Value = 383
$ws.Range("J131").Value = 10102732
$ws.Range("K131").Value = 1149
$ws.Range("L131").Value = 30308196
$ws.Range("M131").Value = 3891
$ws.Range("N131").Value = -30318276
$ws.Range("H141").Value = 4803.5454
$ws.Range("I141").Value = 5334.2856
$ws.Range("J141").Value = 3874.75
$ws.Range("K141").Value = 16002.8568
$ws.Range("L141").Value = 11624.25
$ws.Range("M141").Value = -10822.8568
$ws.Range("N141").Value = -21984.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H97").Value = 1140.4286
$ws.Range("I97").Value = 955.6429000000001
$ws.Range("J97").Value = 1510
$ws.Range("K97").Value = 955.6429000000001
$ws.Range("L97").Value = 1510
$ws.Range("M97").Value = -459.6429000000001
$ws.Range("N97").Value = -2502
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H122").Value = 2966.6667
$ws.Range("I122").Value = 2760
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8280
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5830
$ws.Range("N122").Value = -16900

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 32518.5
$ws.Range("I41").Value = 50000
$ws.Range("J41").Value = 15037
$ws.Range("K41").Value = 50000
$ws.Range("L41").Value = 15037
$ws.Range("M41").Value = -49562
$ws.Range("N41").Value = -15913
$ws.Range("H42").Value = 32500
$ws.Range("I42").Value = 50000
$ws.Range("J42").Value = 15000
$ws.Range("K42").Value = 50000
$ws.Range("L42").Value = 15000
$ws.Range("M42").Value = -49437
$ws.Range("N42").Value = -16126
$ws.Range("H49").Value = 32500
$ws.Range("I49").Value = 50000
$ws.Range("J49").Value = 15000
$ws.Range("K49").Value = 50000
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = -49853
$ws.Range("N49").Value = -15294
$ws.Range("H130").Value = 17848.143
$ws.Range("J130").Value = 17848.143
$ws.Range("L130").Value = 17848.143
$ws.Range("N130").Value = -27888.143
$ws.Range("H132").Value = 2650.0222
$ws.Range("I132").Value = 1972.3422
$ws.Range("J132").Value = 6328.857
$ws.Range("K132").Value = 5917.0266
$ws.Range("L132").Value = 18986.571
$ws.Range("M132").Value = -3387.0266
$ws.Range("N132").Value = -24046.571
$ws.Range("H141").Value = 52915
$ws.Range("J141").Value = 52915
$ws.Range("L141").Value = 52915
$ws.Range("N141").Value = -63275

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 647.3333
$ws.Range("I96").Value = 290
$ws.Range("J96").Value = 826
$ws.Range("K96").Value = 290
$ws.Range("L96").Value = 826
$ws.Range("M96").Value = 1083
$ws.Range("N96").Value = -3572
$ws.Range("H136").Value = 1599.1923
$ws.Range("I136").Value = 753.95
$ws.Range("J136").Value = 4416.6665
$ws.Range("K136").Value = 2261.85
$ws.Range("L136").Value = 13249.9995
$ws.Range("M136").Value = 288.1499999999996
$ws.Range("N136").Value = -18349.9995
